$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModBus")

# ---------------------------------------------------------------------------
# 1) Fix up the cell STYLES first (format-only copy/paste), before touching
#    any values, so that the style-source cells are still carrying their
#    original formatting when we read from them.
#
#    D3:D4 currently carry style 12, D5:D6 style 13, D7:D8 style 11.
#    After the edit they need to become: D3:D4 -> 11, D5:D6 -> 12, D7:D8 -> 13.
#
#    Order matters because the sources and destinations overlap:
#      a) copy style 12 (still sitting on D3) onto D5:D6
#      b) copy style 11 (still sitting on D7) onto D3:D4
#      c) copy style 13 from the untouched B15 cell onto D7:D8
# ---------------------------------------------------------------------------
$ws.Range("D3").Copy()
$ws.Range("D5:D6").PasteSpecial(-4122)

$ws.Range("D7").Copy()
$ws.Range("D3:D4").PasteSpecial(-4122)

$ws.Range("B15").Copy()
$ws.Range("D7:D8").PasteSpecial(-4122)

# F10:G10 go from style 4 to style 1 (copy format from the untouched F11:G11)
$ws.Range("F11:G11").Copy()
$ws.Range("F10:G10").PasteSpecial(-4122)

# F16:G16 go from style 1 to style 4 (copy format from the untouched F9:G9)
$ws.Range("F9:G9").Copy()
$ws.Range("F16:G16").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Now overwrite the cell VALUES to their final contents.
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = "Red Processer"
$ws.Range("E3").Value = "RIO-6/7"
$ws.Range("F3").Value = "Counter, A-stop, & E-stop Reset"

$ws.Range("D4").Value = "Blue Processer"
$ws.Range("E4").Value = "BIO-6/7"
$ws.Range("F4").Value = "Stack Light Green"
$ws.Range("G4").Value = "L0.4"

$ws.Range("D5").ClearContents()
$ws.Range("F5").Value = "Stack Light Orange"
$ws.Range("G5").Value = "L0.3"

$ws.Range("D6").ClearContents()
$ws.Range("F6").Value = "Stack Light Red"
$ws.Range("G6").Value = "L0.2"

$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").Value = "Stack Light Blue"
$ws.Range("G7").Value = "L0.1"

$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").Value = "Stack Buzzer"
$ws.Range("G8").Value = "L0.5"

$ws.Range("F9").Value = "Field Reset Green Light"
$ws.Range("G9").Value = "DS-0/1-GS1"

$ws.Range("F10").Value = "Red Truss Light Outer"
$ws.Range("G10").Value = "BIO-4/5"

$ws.Range("F11").Value = "Red Truss Light Middle"

$ws.Range("F12").Value = "Red Truss Light Inner"
$ws.Range("G12").Value = "BIO-0/1"

$ws.Range("F13").Value = "Blue Truss Light Outer"
$ws.Range("G13").Value = "RIO-4/5"

$ws.Range("F14").Value = "Blue Truss Light Middle"

$ws.Range("C15").Value = "RS-GE-2"
$ws.Range("F15").Value = "Blue Truss Light Inner"
$ws.Range("G15").Value = "RIO-0/1"

$ws.Range("C16").Value = "RS-GE-3"
$ws.Range("F16").ClearContents()
$ws.Range("G16").ClearContents()

$ws.Range("C17").Value = "RS-GE-4"
$ws.Range("C18").Value = "BS-GE-2"
$ws.Range("C19").Value = "BS-GE-3"
$ws.Range("C20").Value = "BS-GE-4"

# ---------------------------------------------------------------------------
# 3) Move the active selection to F3 (was F11 before the edit).
# ---------------------------------------------------------------------------
$ws.Range("F3").Select()
